$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G2").Value = 1.75219
$ws.Range("H2").Value = 5.25657
$ws.Range("I2").Value = 0.1346950717404098
$ws.Range("J2").Value = 0.1418982611249563
$ws.Range("M2").Value = 11.53121
$ws.Range("N2").Value = 34.59363
$ws.Range("O2").Value = 0.04670311854310271
$ws.Range("P2").Value = 0.04697417460197403
$ws.Range("Q2").Value = 20.2048708499
$ws.Range("R2").Value = 181.8438376491
$ws.Range("S2").Value = 0.006290679902664085
$ws.Range("T2").Value = 0.0066655536938002
$ws.Range("G3").Value = 1.75219
$ws.Range("H3").Value = 5.25657
$ws.Range("I3").Value = 0.1346950717404098
$ws.Range("J3").Value = 0.1418982611249563
$ws.Range("O3").Value = 0.04941252241252041
$ws.Range("P3").Value = 0.04969930333854504
$ws.Range("Q3").Value = 21.37702288106
$ws.Range("R3").Value = 192.39320592954
$ws.Range("S3").Value = 0.006655623251229047
$ws.Range("T3").Value = 0.007052244722861274
$ws.Range("G4").Value = 1.75219
$ws.Range("H4").Value = 5.25657
$ws.Range("I4").Value = 0.1346950717404098
$ws.Range("J4").Value = 0.1418982611249563
$ws.Range("M4").Value = 113.969907
$ws.Range("N4").Value = 341.909721
$ws.Range("O4").Value = 0.4615951038067463
$ws.Range("P4").Value = 0.4642741144067919
$ws.Range("Q4").Value = 199.69693134633
$ws.Range("R4").Value = 1797.27238211697
$ws.Range("S4").Value = 0.06217458562227162
$ws.Range("T4").Value = 0.06587968951965277
$ws.Range("G5").Value = 1.75219
$ws.Range("H5").Value = 5.25657
$ws.Range("I5").Value = 0.1346950717404098
$ws.Range("J5").Value = 0.1418982611249563
$ws.Range("M5").Value = 4.2741545
$ws.Range("N5").Value = 8.548309
$ws.Range("O5").Value = 0.01731096253429049
$ws.Range("P5").Value = 0.01160762138918714
$ws.Range("Q5").Value = 7.489130773354999
$ws.Range("R5").Value = 44.93478464013
$ws.Range("S5").Value = 0.002331701340451805
$ws.Range("T5").Value = 0.001647101290922505
$ws.Range("G6").Value = 1.75219
$ws.Range("H6").Value = 5.25657
$ws.Range("I6").Value = 0.1346950717404098
$ws.Range("J6").Value = 0.1418982611249563
$ws.Range("M6").Value = 104.9290516666667
$ws.Range("N6").Value = 314.787155
$ws.Range("O6").Value = 0.4249782927033401
$ws.Range("P6").Value = 0.4274447862635018
$ws.Range("Q6").Value = 183.8556350398167
$ws.Range("R6").Value = 1654.70071535835
$ws.Range("S6").Value = 0.05724248162379329
$ws.Range("T6").Value = 0.0606536718977195
$ws.Range("I7").Value = 0.1557790731975008
$ws.Range("J7").Value = 0.1641097875428132
$ws.Range("M7").Value = 11.53121
$ws.Range("N7").Value = 34.59363
$ws.Range("O7").Value = 0.04670311854310271
$ws.Range("P7").Value = 0.04697417460197403
$ws.Range("Q7").Value = 23.36756656649333
$ws.Range("R7").Value = 210.30809909844
$ws.Range("S7").Value = 0.007275368522077553
$ws.Range("T7").Value = 0.00770892181392897
$ws.Range("I8").Value = 0.1557790731975008
$ws.Range("J8").Value = 0.1641097875428132
$ws.Range("O8").Value = 0.04941252241252041
$ws.Range("P8").Value = 0.04969930333854504
$ws.Range("S8").Value = 0.007697436945773166
$ws.Range("T8").Value = 0.008156142111914453
$ws.Range("I9").Value = 0.1557790731975008
$ws.Range("J9").Value = 0.1641097875428132
$ws.Range("M9").Value = 113.969907
$ws.Range("N9").Value = 341.909721
$ws.Range("O9").Value = 0.4615951038067463
$ws.Range("P9").Value = 0.4642741144067919
$ws.Range("Q9").Value = 230.955761658972
$ws.Range("R9").Value = 2078.601854930748
$ws.Range("S9").Value = 0.0719068574635191
$ws.Range("T9").Value = 0.07619192627692636
$ws.Range("I10").Value = 0.1557790731975008
$ws.Range("J10").Value = 0.1641097875428132
$ws.Range("M10").Value = 4.2741545
$ws.Range("N10").Value = 8.548309
$ws.Range("O10").Value = 0.01731096253429049
$ws.Range("P10").Value = 0.01160762138918714
$ws.Range("Q10").Value = 8.661414525815333
$ws.Range("R10").Value = 51.968487154892
$ws.Range("S10").Value = 0.002696685699748432
$ws.Range("T10").Value = 0.001904924280056916
$ws.Range("I11").Value = 0.1557790731975008
$ws.Range("J11").Value = 0.1641097875428132
$ws.Range("M11").Value = 104.9290516666667
$ws.Range("N11").Value = 314.787155
$ws.Range("O11").Value = 0.4249782927033401
$ws.Range("P11").Value = 0.4274447862635018
$ws.Range("Q11").Value = 212.6348058512378
$ws.Range("R11").Value = 1913.71325266114
$ws.Range("S11").Value = 0.06620272456638253
$ws.Range("T11").Value = 0.07014787305998647
$ws.Range("G12").Value = 4.488144
$ws.Range("H12").Value = 13.464432
$ws.Range("I12").Value = 0.3450144550883694
$ws.Range("J12").Value = 0.3634650518941472
$ws.Range("M12").Value = 11.53121
$ws.Range("N12").Value = 34.59363
$ws.Range("O12").Value = 0.04670311854310271
$ws.Range("P12").Value = 0.04697417460197403
$ws.Range("Q12").Value = 51.75373097424
$ws.Range("R12").Value = 465.78357876816
$ws.Range("S12").Value = 0.0161132509950761
$ws.Range("T12").Value = 0.01707347080939122
$ws.Range("G13").Value = 4.488144
$ws.Range("H13").Value = 13.464432
$ws.Range("I13").Value = 0.3450144550883694
$ws.Range("J13").Value = 0.3634650518941472
$ws.Range("O13").Value = 0.04941252241252041
$ws.Range("P13").Value = 0.04969930333854504
$ws.Range("Q13").Value = 54.75613773705599
$ws.Range("R13").Value = 492.805239633504
$ws.Range("S13").Value = 0.01704803449469757
$ws.Range("T13").Value = 0.01806395986704723
$ws.Range("G14").Value = 4.488144
$ws.Range("H14").Value = 13.464432
$ws.Range("I14").Value = 0.3450144550883694
$ws.Range("J14").Value = 0.3634650518941472
$ws.Range("M14").Value = 113.969907
$ws.Range("N14").Value = 341.909721
$ws.Range("O14").Value = 0.4615951038067463
$ws.Range("P14").Value = 0.4642741144067919
$ws.Range("Q14").Value = 511.513354282608
$ws.Range("R14").Value = 4603.620188543472
$ws.Range("S14").Value = 0.1592569832113439
$ws.Range("T14").Value = 0.1687474150859739
$ws.Range("G15").Value = 4.488144
$ws.Range("H15").Value = 13.464432
$ws.Range("I15").Value = 0.3450144550883694
$ws.Range("J15").Value = 0.3634650518941472
$ws.Range("M15").Value = 4.2741545
$ws.Range("N15").Value = 8.548309
$ws.Range("O15").Value = 0.01731096253429049
$ws.Range("P15").Value = 0.01160762138918714
$ws.Range("Q15").Value = 19.183020874248
$ws.Range("R15").Value = 115.098125245488
$ws.Range("S15").Value = 0.005972532305823413
$ws.Range("T15").Value = 0.004218964710588518
$ws.Range("G16").Value = 4.488144
$ws.Range("H16").Value = 13.464432
$ws.Range("I16").Value = 0.3450144550883694
$ws.Range("J16").Value = 0.3634650518941472
$ws.Range("M16").Value = 104.9290516666667
$ws.Range("N16").Value = 314.787155
$ws.Range("O16").Value = 0.4249782927033401
$ws.Range("P16").Value = 0.4274447862635018
$ws.Range("Q16").Value = 470.93669366344
$ws.Range("R16").Value = 4238.43024297096
$ws.Range("S16").Value = 0.1466236540814284
$ws.Range("T16").Value = 0.1553612414211463
$ws.Range("G17").Value = 1.981064
$ws.Range("H17").Value = 3.962128
$ws.Range("I17").Value = 0.152289168185153
$ws.Range("J17").Value = 0.1069555001749241
$ws.Range("M17").Value = 11.53121
$ws.Range("N17").Value = 34.59363
$ws.Range("O17").Value = 0.04670311854310271
$ws.Range("P17").Value = 0.04697417460197403
$ws.Range("Q17").Value = 22.84406500744
$ws.Range("R17").Value = 137.06439004464
$ws.Range("S17").Value = 0.007112379074581707
$ws.Range("T17").Value = 0.005024146339858348
$ws.Range("G18").Value = 1.981064
$ws.Range("H18").Value = 3.962128
$ws.Range("I18").Value = 0.152289168185153
$ws.Range("J18").Value = 0.1069555001749241
$ws.Range("O18").Value = 0.04941252241252041
$ws.Range("P18").Value = 0.04969930333854504
$ws.Range("Q18").Value = 24.169325505136
$ws.Range("R18").Value = 145.015953030816
$ws.Range("S18").Value = 0.007524991936132964
$ws.Range("T18").Value = 0.005315613846919359
$ws.Range("G19").Value = 1.981064
$ws.Range("H19").Value = 3.962128
$ws.Range("I19").Value = 0.152289168185153
$ws.Range("J19").Value = 0.1069555001749241
$ws.Range("M19").Value = 113.969907
$ws.Range("N19").Value = 341.909721
$ws.Range("O19").Value = 0.4615951038067463
$ws.Range("P19").Value = 0.4642741144067919
$ws.Range("Q19").Value = 225.781679841048
$ws.Range("R19").Value = 1354.690079046288
$ws.Range("S19").Value = 0.07029593439706874
$ws.Range("T19").Value = 0.04965667012464835
$ws.Range("G20").Value = 1.981064
$ws.Range("H20").Value = 3.962128
$ws.Range("I20").Value = 0.152289168185153
$ws.Range("J20").Value = 0.1069555001749241
$ws.Range("M20").Value = 4.2741545
$ws.Range("N20").Value = 8.548309
$ws.Range("O20").Value = 0.01731096253429049
$ws.Range("P20").Value = 0.01160762138918714
$ws.Range("Q20").Value = 8.467373610388
$ws.Range("R20").Value = 33.869494441552
$ws.Range("S20").Value = 0.002636272084831447
$ws.Range("T20").Value = 0.001241498951521658
$ws.Range("G21").Value = 1.981064
$ws.Range("H21").Value = 3.962128
$ws.Range("I21").Value = 0.152289168185153
$ws.Range("J21").Value = 0.1069555001749241
$ws.Range("M21").Value = 104.9290516666667
$ws.Range("N21").Value = 314.787155
$ws.Range("O21").Value = 0.4249782927033401
$ws.Range("P21").Value = 0.4274447862635018
$ws.Range("Q21").Value = 207.8711668109733
$ws.Range("R21").Value = 1247.22700086584
$ws.Range("S21").Value = 0.06471959069253815
$ws.Range("T21").Value = 0.04571757091197635
$ws.Range("G22").Value = 2.760707333333334
$ws.Range("H22").Value = 8.282122000000001
$ws.Range("I22").Value = 0.2122222317885668
$ws.Range("J22").Value = 0.2235713992631593
$ws.Range("M22").Value = 11.53121
$ws.Range("N22").Value = 34.59363
$ws.Range("O22").Value = 0.04670311854310271
$ws.Range("P22").Value = 0.04697417460197403
$ws.Range("Q22").Value = 31.83429600920667
$ws.Range("R22").Value = 286.50866408286
$ws.Range("S22").Value = 0.009911440048703257
$ws.Range("T22").Value = 0.01050208194499529
$ws.Range("G23").Value = 2.760707333333334
$ws.Range("H23").Value = 8.282122000000001
$ws.Range("I23").Value = 0.2122222317885668
$ws.Range("J23").Value = 0.2235713992631593
$ws.Range("O23").Value = 0.04941252241252041
$ws.Range("P23").Value = 0.04969930333854504
$ws.Range("Q23").Value = 33.68110982974267
$ws.Range("R23").Value = 303.129988467684
$ws.Range("S23").Value = 0.01048643578468766
$ws.Range("T23").Value = 0.01111134278980272
$ws.Range("G24").Value = 2.760707333333334
$ws.Range("H24").Value = 8.282122000000001
$ws.Range("I24").Value = 0.2122222317885668
$ws.Range("J24").Value = 0.2235713992631593
$ws.Range("M24").Value = 113.969907
$ws.Range("N24").Value = 341.909721
$ws.Range("O24").Value = 0.4615951038067463
$ws.Range("P24").Value = 0.4642741144067919
$ws.Range("Q24").Value = 314.637558034218
$ws.Range("R24").Value = 2831.738022307962
$ws.Range("S24").Value = 0.09796074311254288
$ws.Range("T24").Value = 0.1037984133995906
$ws.Range("G25").Value = 2.760707333333334
$ws.Range("H25").Value = 8.282122000000001
$ws.Range("I25").Value = 0.2122222317885668
$ws.Range("J25").Value = 0.2235713992631593
$ws.Range("M25").Value = 4.2741545
$ws.Range("N25").Value = 8.548309
$ws.Range("O25").Value = 0.01731096253429049
$ws.Range("P25").Value = 0.01160762138918714
$ws.Range("Q25").Value = 11.79968967194967
$ws.Range("R25").Value = 70.798138031698
$ws.Range("S25").Value = 0.003673771103435393
$ws.Range("T25").Value = 0.002595132156097547
$ws.Range("G26").Value = 2.760707333333334
$ws.Range("H26").Value = 8.282122000000001
$ws.Range("I26").Value = 0.2122222317885668
$ws.Range("J26").Value = 0.2235713992631593
$ws.Range("M26").Value = 104.9290516666667
$ws.Range("N26").Value = 314.787155
$ws.Range("O26").Value = 0.4249782927033401
$ws.Range("P26").Value = 0.4274447862635018
$ws.Range("Q26").Value = 289.678402415879
$ws.Range("R26").Value = 2607.10562174291
$ws.Range("S26").Value = 0.09018984173919764
$ws.Range("T26").Value = 0.09556442897267314
Write-Output "Updated 278 cells"
